# Updates cryptocurrency price/volume data on the Coin sheet,
# mirroring the latest scrape (including two coin re-rankings
# where SuiNetwork/EthereumClassic and Filecoin/Bittensor swapped rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.928.05'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.08%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.344.05'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.68%  '
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '540.86'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.16'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.31%  '
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.563'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.85%  '
$ws.Range("E9").Value = '  +0.72%  '
$ws.Range("E10").Value = '  +2.02%  '
$ws.Range("E12").Value = '  +0.54%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.78'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.15%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.759.37'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.23%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '57.869.14'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.00%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000135'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.66%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.349.96'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.60%  '
$ws.Range("E18").Value = '  +1.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.30'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '329.15'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.70%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.74'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.78%  '
$ws.Range("E22").Value = '  +0.30%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '63.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.48%  '
$ws.Range("E24").Value = '  -3.31%  '
$ws.Range("E25").Value = '  +0.21%  '
$ws.Range("E26").Value = '  -2.02%  '
$ws.Range("E27").Value = '  -5.33%  '
$ws.Range("E28").Value = '  +0.25%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '170.18'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.26%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0735'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.23%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.12'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.75%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.32'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.29%  '
$ws.Range("B33").Value = 'SuiNetwork'
$ws.Range("C33").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.01'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.87%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("E36").Value = '  +1.32%  '
$ws.Range("E37").Value = '  -1.99%  '
$ws.Range("E38").Value = '  -0.55%  '
$ws.Range("E39").Value = '  -0.81%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '141.23'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.69%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.377'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.45%  '
$ws.Range("B42").Value = 'Bittensor'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '288.74'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.84%  '
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.63'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.34%  '
$ws.Range("E44").Value = '  +1.40%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0509'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.74%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.04'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.27%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.566'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.68%  '
$ws.Range("E48").Value = '  +1.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.381'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.30%  '
$ws.Range("E50").Value = '  +0.30%  '
$ws.Range("E51").Value = '  +0.95%  '
